# Add "arrival_location_url" as a new column to the Clients import table,
# and fix the misspelled "chronic_disease_discription" header to
# "chronic_disease_description".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Clients"
$ws2 = $wb.Worksheets.Item(2)   # "Relatives"

# The Clients sheet's table currently spans A1:M1047180 (13 columns).
# Grow it by one column so a 14th column (N) is included.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:N1047179"))

# Give the new column N a header first (so the shared-string table gets
# "arrival_location_url" registered before the renamed column's string),
# then fix column M's header text. Setting the header cell values updates
# both the worksheet cells and the table's column names.
$ws1.Range("N1").Value = "arrival_location_url"
$ws1.Range("M1").Value = "chronic_disease_description"

# Match the column width used for the new column (closest value reachable
# through the ColumnWidth/character-unit rounding the host performs).
$ws1.Columns.Item(14).ColumnWidth = 22.67

# Update the view state: Clients becomes the active/selected tab with M6
# selected, Relatives is no longer the active tab and has I9 selected.
[void]$ws2.Activate()
[void]$ws2.Range("I9").Select()
[void]$ws1.Activate()
[void]$ws1.Range("M6").Select()
